$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save the original values of rows 26, 27, 28 for the columns that change
# (A, B, E, F, G, H, Q, R) before overwriting anything.
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")
$rows = @(26, 27, 28)

$orig = @{}
foreach ($r in $rows) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Cyclic rotation observed in the diff:
#   new row 26 <- old row 28
#   new row 27 <- old row 26
#   new row 28 <- old row 27
$mapping = @{ 26 = 28; 27 = 26; 28 = 27 }

foreach ($destRow in $rows) {
    $srcRow = $mapping[$destRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $orig[$srcRow][$c]
    }
}
